# Weekly fruta/hortaliza update: prepend a new week's pair of rows
# (Primera / Segunda) for "Pepino ensalada" - Agrícola del Norte S.A. de Arica,
# pushing the existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 463, shifting
# rows 463:513 down to 465:515 (keeps all existing formatting/styles).
$ws.Rows.Item(463).Resize(2).Insert()

# --- New row 463: "Primera" quality for the new week ---
$ws.Range("A463").Value = 1
$ws.Range("B463").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C463").Value = "Arica y Parinacota"
$ws.Range("D463").Value = 45124
$ws.Range("E463").Value = 15
$ws.Range("F463").Value = 100112043
$ws.Range("G463").Value = "Pepino ensalada"
$ws.Range("H463").Value = "Sin especificar"
$ws.Range("I463").Value = "Primera"
$ws.Range("J463").Value = 340
$ws.Range("K463").Value = 9000
$ws.Range("L463").Value = 10000
$ws.Range("M463").Value = 9441
$ws.Range("N463").Value = "$/caja 70 unidades"
$ws.Range("O463").Value = "Región de Arica y Parinacota"
$ws.Range("P463").Value = 135
$ws.Range("Q463").Value = 70
$ws.Range("R463").Value = "Hortaliza"

# --- New row 464: "Segunda" quality for the new week ---
$ws.Range("A464").Value = 1
$ws.Range("B464").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C464").Value = "Arica y Parinacota"
$ws.Range("D464").Value = 45124
$ws.Range("E464").Value = 15
$ws.Range("F464").Value = 100112043
$ws.Range("G464").Value = "Pepino ensalada"
$ws.Range("H464").Value = "Sin especificar"
$ws.Range("I464").Value = "Segunda"
$ws.Range("J464").Value = 350
$ws.Range("K464").Value = 7000
$ws.Range("L464").Value = 8000
$ws.Range("M464").Value = 7286
$ws.Range("N464").Value = "$/caja 100 unidades"
$ws.Range("O464").Value = "Región de Arica y Parinacota"
$ws.Range("P464").Value = 73
$ws.Range("Q464").Value = 100
$ws.Range("R464").Value = "Hortaliza"
